$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text: replace spaces with underscores
$ws.Range("B1").Value = "Waste_Generated"
$ws.Range("C1").Value = "Waste_Recycled"
$ws.Range("D1").Value = "Waste_Disposed"
$ws.Range("E1").Value = "Overall_Recycling_Rate$([char]0x00A0)"
$ws.Range("F1").Value = "Overall_Recycling_Rate_Without_Construction_&_Demolition_Waste"

# Update row 1 height
$ws.Rows.Item(1).RowHeight = 65.25

# Update the selection
$ws.Range("F2").Select()
